$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Data row (entered first so shared-string table order matches)
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"

# Header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Zoom
$ws.Application.ActiveWindow.Zoom = 205

# Selection
$ws.Range("A3").Select()

# Column width for B (auto-fit-ish best fit)
$ws.Columns.Item(2).ColumnWidth = 11.140625
